$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 2016530.8
$ws.Range("I38").Value = 4608432.5
$ws.Range("J38").Value = 607.2222
$ws.Range("K38").Value = 13825297.5
$ws.Range("L38").Value = 1821.6666
$ws.Range("M38").Value = -13824925.5
$ws.Range("N38").Value = -2565.6666
$ws.Range("H39").Value = 1443190
$ws.Range("I39").Value = 1984242.5
$ws.Range("J39").Value = 383.33334
$ws.Range("K39").Value = 5952727.5
$ws.Range("L39").Value = 1150.00002
$ws.Range("M39").Value = -5952431.5
$ws.Range("N39").Value = -1742.00002
$ws.Range("H40").Value = 1805.05
$ws.Range("I40").Value = 1673.4
$ws.Range("J40").Value = 2200
$ws.Range("K40").Value = 1673.4
$ws.Range("L40").Value = 2200
$ws.Range("M40").Value = -1498.4
$ws.Range("N40").Value = -2550
$ws.Range("H42").Value = 1562924.1
$ws.Range("I42").Value = 5000211.5
$ws.Range("J42").Value = 520.7273
$ws.Range("K42").Value = 15000634.5
$ws.Range("L42").Value = 1562.1819
$ws.Range("M42").Value = -15000404.5
$ws.Range("N42").Value = -2022.1819
$ws.Range("H58").Value = 1378935.5
$ws.Range("I58").Value = 2525548.2
$ws.Range("J58").Value = 3000
$ws.Range("K58").Value = 7576644.600000001
$ws.Range("L58").Value = 9000
$ws.Range("M58").Value = -7576494.600000001
$ws.Range("N58").Value = -9300
$ws.Range("H74").Value = 3811.3333
$ws.Range("I74").Value = 3287.875
$ws.Range("K74").Value = 3287.875
$ws.Range("M74").Value = -2351.875
$ws.Range("H77").Value = 3811.3333
$ws.Range("I77").Value = 3287.875
$ws.Range("K77").Value = 16439.375
$ws.Range("M77").Value = -11759.375
$ws.Range("H100").Value = 2715
$ws.Range("I100").Value = 2876.25
$ws.Range("J100").Value = 2500
$ws.Range("K100").Value = 2876.25
$ws.Range("L100").Value = 2500
$ws.Range("M100").Value = -2335.25
$ws.Range("N100").Value = -3582
$ws.Range("H137").Value = 1430.4894
$ws.Range("I137").Value = 1505.1538
$ws.Range("K137").Value = 4515.4614
$ws.Range("M137").Value = -1965.4614

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 21842.06
$ws.Range("I32").Value = 3498.1487
$ws.Range("K32").Value = 3498.1487
$ws.Range("M32").Value = -3211.1487
$ws.Range("H121").Value = 49980
$ws.Range("J121").Value = 49980
$ws.Range("L121").Value = 49980
$ws.Range("N121").Value = -53474
$ws.Range("H132").Value = 2794.0208
$ws.Range("I132").Value = 3784.3447
$ws.Range("J132").Value = 1282.4736
$ws.Range("K132").Value = 11353.0341
$ws.Range("L132").Value = 3847.4208
$ws.Range("M132").Value = -8823.034100000001
$ws.Range("N132").Value = -8907.4208

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 77203.08
$ws.Range("I94").Value = 77203.08
$ws.Range("K94").Value = 77203.08
$ws.Range("M94").Value = -76752.08
$ws.Range("H107").Value = 16689227
$ws.Range("I107").Value = 20860276
$ws.Range("J107").Value = 5030.75
$ws.Range("K107").Value = 20860276
$ws.Range("L107").Value = 5030.75
$ws.Range("M107").Value = -20858356
$ws.Range("N107").Value = -8870.75

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 17312.8
$ws.Range("I31").Value = 44540.914
$ws.Range("J31").Value = 2402.1667
$ws.Range("K31").Value = 44540.914
$ws.Range("L31").Value = 2402.1667
$ws.Range("M31").Value = -44245.914
$ws.Range("N31").Value = -2992.1667
$ws.Range("H34").Value = 17312.8
$ws.Range("I34").Value = 44540.914
$ws.Range("J34").Value = 2402.1667
$ws.Range("K34").Value = 44540.914
$ws.Range("L34").Value = 2402.1667
$ws.Range("M34").Value = -44338.914
$ws.Range("N34").Value = -2806.1667
$ws.Range("H62").Value = 3705957
$ws.Range("I62").Value = 6946432
$ws.Range("J62").Value = 2557.1428
$ws.Range("K62").Value = 6946432
$ws.Range("L62").Value = 2557.1428
$ws.Range("M62").Value = -6945808
$ws.Range("N62").Value = -3805.1428
$ws.Range("H65").Value = 3705957
$ws.Range("I65").Value = 6946432
$ws.Range("J65").Value = 2557.1428
$ws.Range("K65").Value = 34732160
$ws.Range("L65").Value = 12785.714
$ws.Range("M65").Value = -34729040
$ws.Range("N65").Value = -19025.714
$ws.Range("H108").Value = 39899.5
$ws.Range("J108").Value = 39899.5
$ws.Range("L108").Value = 39899.5
$ws.Range("N108").Value = -47579.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1411.7548
$ws.Range("I131").Value = 1134.2858
$ws.Range("J131").Value = 1453.9783
$ws.Range("K131").Value = 3402.8574
$ws.Range("L131").Value = 4361.9349
$ws.Range("M131").Value = 1637.1426
$ws.Range("N131").Value = -14441.9349
$ws.Range("H134").Value = 4899.615
$ws.Range("I134").Value = 2608.5715
$ws.Range("J134").Value = 7572.5
$ws.Range("K134").Value = 7825.7145
$ws.Range("L134").Value = 22717.5
$ws.Range("M134").Value = -2755.7145
$ws.Range("N134").Value = -32857.5
$ws.Range("H139").Value = 1657.6207
$ws.Range("I139").Value = 1014
$ws.Range("J139").Value = 4747
$ws.Range("K139").Value = 3042
$ws.Range("L139").Value = 14241
$ws.Range("M139").Value = 2098
$ws.Range("N139").Value = -24521

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 58825364
$ws.Range("I97").Value = 90911110
$ws.Range("K97").Value = 90911110
$ws.Range("M97").Value = -90910614
$ws.Range("H102").Value = 355185.53
$ws.Range("I102").Value = 2410.3333
$ws.Range("J102").Value = 3000999.5
$ws.Range("K102").Value = 2410.3333
$ws.Range("L102").Value = 3000999.5
$ws.Range("M102").Value = -788.3332999999998
$ws.Range("N102").Value = -3004243.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 947701
$ws.Range("I55").Value = 1623648
$ws.Range("K55").Value = 1623648
$ws.Range("M55").Value = -1623475
$ws.Range("H68").Value = 2338.4092
$ws.Range("I68").Value = 1395.9286
$ws.Range("J68").Value = 3987.75
$ws.Range("K68").Value = 1395.9286
$ws.Range("L68").Value = 3987.75
$ws.Range("M68").Value = -646.9286
$ws.Range("N68").Value = -5485.75
$ws.Range("H71").Value = 2338.4092
$ws.Range("I71").Value = 1395.9286
$ws.Range("J71").Value = 3987.75
$ws.Range("K71").Value = 6979.643
$ws.Range("L71").Value = 19938.75
$ws.Range("M71").Value = -3235.643
$ws.Range("N71").Value = -27426.75

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 1404795.5
$ws.Range("I62").Value = 2900490.5
$ws.Range("J62").Value = 2581.25
$ws.Range("K62").Value = 2900490.5
$ws.Range("L62").Value = 2581.25
$ws.Range("M62").Value = -2899866.5
$ws.Range("N62").Value = -3829.25
$ws.Range("H65").Value = 1404795.5
$ws.Range("I65").Value = 2900490.5
$ws.Range("J65").Value = 2581.25
$ws.Range("K65").Value = 14502452.5
$ws.Range("L65").Value = 12906.25
$ws.Range("M65").Value = -14499332.5
$ws.Range("N65").Value = -19146.25
$ws.Range("H107").Value = 37148.848
$ws.Range("I107").Value = 9912.714
$ws.Range("J107").Value = 84812.086
$ws.Range("K107").Value = 29738.142
$ws.Range("L107").Value = 254436.258
$ws.Range("M107").Value = -27818.142
$ws.Range("N107").Value = -258276.258

Write-Output "Applied all updates"